# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must stay text (prices look numeric, e.g. "238.77")
# Excel's COM .Value setter auto-coerces plain numeric-looking strings into
# numbers, so we force the number format to Text first, then reset the
# style back to Normal afterward so we don't leave a stray number format
# applied to the cell (matching the source which uses default styling).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "97.345.00"
$ws.Range("E2").Value = "  +0.52%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.736.38"
$ws.Range("E3").Value = "  +1.83%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - Solana
Set-TextValue $ws.Range("D5") "238.77"
$ws.Range("E5").Value = "  -0.58%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.53%  "

# Row 7 - BNB
Set-TextValue $ws.Range("D7") "662.97"
$ws.Range("E7").Value = "  +0.74%  "

# Row 8 - Dogecoin
Set-TextValue $ws.Range("D8") "0.441"
$ws.Range("E8").Value = "  +2.96%  "

# Row 9 - was Cardano, now USDC
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D9") "0.999"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10 - was USDC, now Cardano
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D10") "1.07"
$ws.Range("E10").Value = "  -2.04%  "

# Row 11 - LidoStakedEther
Set-TextValue $ws.Range("D11") "3.735.94"
$ws.Range("E11").Value = "  +1.88%  "

# Row 12 - ShibaInu
Set-TextValue $ws.Range("D12") "0.0000321"
$ws.Range("E12").Value = "  +17.95%  "

# Row 13 - Avalanche
Set-TextValue $ws.Range("D13") "45.11"
$ws.Range("E13").Value = "  -0.94%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.46%  "

# Row 15 - Toncoin
Set-TextValue $ws.Range("D15") "6.93"
$ws.Range("E15").Value = "  +2.31%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D16") "4.433.28"
$ws.Range("E16").Value = "  +1.83%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "97.023.67"
$ws.Range("E17").Value = "  +0.47%  "

# Row 18 - Polkadot
Set-TextValue $ws.Range("D18") "9.07"
$ws.Range("E18").Value = "  +2.02%  "

# Row 19 - WrappedEther
Set-TextValue $ws.Range("D19") "3.738.54"
$ws.Range("E19").Value = "  +1.78%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "13.13"
$ws.Range("E20").Value = "  +2.82%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("D21") "18.86"
$ws.Range("E21").Value = "  +0.54%  "

# Row 22 - Stellar
Set-TextValue $ws.Range("D22") "0.506"
$ws.Range("E22").Value = "  -4.12%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "530.36"
$ws.Range("E23").Value = "  -0.64%  "

# Row 25 - PEPE
Set-TextValue $ws.Range("D25") "0.0000230"
$ws.Range("E25").Value = "  +11.49%  "

# Row 26 - NEARProtocol
Set-TextValue $ws.Range("D26") "6.91"
$ws.Range("E26").Value = "  -3.43%  "

# Row 27 - Litecoin
Set-TextValue $ws.Range("D27") "107.43"
$ws.Range("E27").Value = "  +4.93%  "

# Row 28 - was Hedera, now Aptos
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D28") "13.67"
$ws.Range("E28").Value = "  +1.61%  "

# Row 29 - was Aptos, now Hedera
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D29") "0.190"
$ws.Range("E29").Value = "  +13.37%  "

# Row 30 - WrappedeETH
Set-TextValue $ws.Range("D30") "3.932.07"
$ws.Range("E30").Value = "  +1.77%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D31") "12.84"
$ws.Range("E31").Value = "  +3.27%  "

# Row 32 - PancakeSwap
Set-TextValue $ws.Range("D32") "3.05"
$ws.Range("E32").Value = "  +0.29%  "

# Row 33 - Dai
Set-TextValue $ws.Range("D33") "1.00"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34 - Cronos
$ws.Range("E34").Value = "  +3.25%  "

# Row 35 - Fetch.AI
$ws.Range("E35").Value = "  -3.71%  "

# Row 36 - EthereumClassic
Set-TextValue $ws.Range("D36") "32.59"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D37") "0.998"
$ws.Range("E37").Value = "  -0.20%  "

# Row 38 - Bittensor
Set-TextValue $ws.Range("D38") "642.34"
$ws.Range("E38").Value = "  -3.38%  "

# Row 39 - PolygonEcosystemToken
Set-TextValue $ws.Range("D39") "0.595"
$ws.Range("E39").Value = "  -0.19%  "

# Row 40 - RenderToken
Set-TextValue $ws.Range("D40") "8.80"
$ws.Range("E40").Value = "  -0.53%  "

# Row 42 - Kaspa
Set-TextValue $ws.Range("D42") "0.167"
$ws.Range("E42").Value = "  +4.20%  "

# Row 43 - was Filecoin, now ImmutableX
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D43") "2.04"
$ws.Range("E43").Value = "  +2.14%  "

# Row 44 - was ImmutableX, now Filecoin
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D44") "6.79"
$ws.Range("E44").Value = "  +3.69%  "

# Row 45 - EnergySwap
Set-TextValue $ws.Range("D45") "40.68"
$ws.Range("E45").Value = "  +4.81%  "

# Row 46 - ARBITRUM
Set-TextValue $ws.Range("D46") "0.978"
$ws.Range("E46").Value = "  +1.98%  "

# Row 47 - Algorand
Set-TextValue $ws.Range("D47") "0.477"
$ws.Range("E47").Value = "  +9.41%  "

# Row 48 - VeChain
Set-TextValue $ws.Range("D48") "0.0458"
$ws.Range("E48").Value = "  -0.79%  "

# Row 49 - Stacks
Set-TextValue $ws.Range("D49") "2.39"
$ws.Range("E49").Value = "  +2.01%  "

# Row 50 - WhiteBITCoin
Set-TextValue $ws.Range("D50") "23.64"

# Row 51 - Cosmos
Set-TextValue $ws.Range("D51") "8.70"
$ws.Range("E51").Value = "  -0.13%  "
